{"js": "const replacements = [\n  [\"2025-12-09 Tuesday\", \"2025-12-10 Wednesday\"],\n  [\"32\u00d748=1536\", \"14\u00d782=1148\"],\n  [\"84\u00d763=5292\", \"55\u00d719=1045\"],\n  [\"61\u00d794=5734\", \"93\u00d796=8928\"],\n  [\"24\u00d782=1968\", \"50\u00d741=2050\"],\n  [\"77\u00d758=4466\", \"13\u00d773=949\"],\n  [\"22\u00d794=2068\", \"64\u00d717=1088\"],\n  [\"92\u00d767=6164\", \"63\u00d732=2016\"],\n  [\"30\u00d736=1080\", \"25\u00d782=2050\"],\n  [\"67\u00d762=4154\", \"48\u00d747=2256\"],\n  [\"16\u00d742=672\", \"56\u00d732=1792\"],\n  [\"82\u00d760=4920\", \"88\u00d726=2288\"],\n  [\"58\u00d796=5568\", \"30\u00d712=360\"],\n  [\"53\u00d772=3816\", \"94\u00d781=7614\"],\n  [\"22\u00d777=1694\", \"65\u00d746=2990\"],\n  [\"53\u00d793=4929\", \"35\u00d754=1890\"],\n  [\"93\u00d728=2604\", \"73\u00d740=2920\"],\n  [\"38\u00d718=684\", \"99\u00d787=8613\"],\n  [\"64\u00d781=5184\", \"30\u00d792=2760\"],\n  [\"91\u00d782=7462\", \"36\u00d778=2808\"],\n  [\"25\u00d787=2175\", \"50\u00d768=3400\"],\n  [\"20\u00d743=860\", \"64\u00d763=4032\"],\n  [\"36\u00d733=1188\", \"95\u00d722=2090\"],\n  [\"45\u00d731=1395\", \"21\u00d765=1365\"],\n  [\"74\u00d772=5328\", \"64\u00d739=2496\"],\n  [\"15\u00d760=900\", \"25\u00d758=1450\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load('items');\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}", "ps1": "$d = $word.ActiveDocument\n$wdReplaceAll = 2\n$wdFindContinue = 1\n\n$pairs = @(\n    ,@(\"2025-12-09 Tuesday\", \"2025-12-10 Wednesday\")\n    ,@(\"32\u00d748=1536\", \"14\u00d782=1148\")\n    ,@(\"84\u00d763=5292\", \"55\u00d719=1045\")\n    ,@(\"61\u00d794=5734\", \"93\u00d796=8928\")\n    ,@(\"24\u00d782=1968\", \"50\u00d741=2050\")\n    ,@(\"77\u00d758=4466\", \"13\u00d773=949\")\n    ,@(\"22\u00d794=2068\", \"64\u00d717=1088\")\n    ,@(\"92\u00d767=6164\", \"63\u00d732=2016\")\n    ,@(\"30\u00d736=1080\", \"25\u00d782=2050\")\n    ,@(\"67\u00d762=4154\", \"48\u00d747=2256\")\n    ,@(\"16\u00d742=672\", \"56\u00d732=1792\")\n    ,@(\"82\u00d760=4920\", \"88\u00d726=2288\")\n    ,@(\"58\u00d796=5568\", \"30\u00d712=360\")\n    ,@(\"53\u00d772=3816\", \"94\u00d781=7614\")\n    ,@(\"22\u00d777=1694\", \"65\u00d746=2990\")\n    ,@(\"53\u00d793=4929\", \"35\u00d754=1890\")\n    ,@(\"93\u00d728=2604\", \"73\u00d740=2920\")\n    ,@(\"38\u00d718=684\", \"99\u00d787=8613\")\n    ,@(\"64\u00d781=5184\", \"30\u00d792=2760\")\n    ,@(\"91\u00d782=7462\", \"36\u00d778=2808\")\n    ,@(\"25\u00d787=2175\", \"50\u00d768=3400\")\n    ,@(\"20\u00d743=860\", \"64\u00d763=4032\")\n    ,@(\"36\u00d733=1188\", \"95\u00d722=2090\")\n    ,@(\"45\u00d731=1395\", \"21\u00d765=1365\")\n    ,@(\"74\u00d772=5328\", \"64\u00d739=2496\")\n    ,@(\"15\u00d760=900\", \"25\u00d758=1450\")\n)\n\nforeach ($pair in $pairs) {\n    $findText = $pair[0]\n    $replaceText = $pair[1]\n    $r = $d.Content\n    $r.Find.ClearFormatting()\n    $r.Find.Replacement.ClearFormatting()\n    $r.Find.Text = $findText\n    $r.Find.Replacement.Text = $replaceText\n    $r.Find.Execute($findText, $true, $true, $false, $false, $false, $true, $wdFindContinue, $false, $replaceText, $wdReplaceAll)\n}\n"}
